$p = $ppt.ActivePresentation
try {
  $r = $p.Designs.Load("ppt/theme/theme1.xml")
  Write-Output ("Load result: " + $r)
} catch {
  Write-Output ("Load error: " + $_.Exception.Message)
}
